$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - unique values
$ws.Range("B2").Value = 5
$ws.Range("D2").Value = 0.14
$ws.Range("F2").Value = 0.26
$ws.Range("H2").Value = 0.15

# Rows 3-18 - shared values
for ($r = 3; $r -le 18; $r++) {
    $ws.Range("B$r").Value = 0.3125
    $ws.Range("D$r").Value = 0.05
    $ws.Range("F$r").Value = 0.04
    $ws.Range("H$r").Value = 0.02
}
